$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "date header" / weekly tracking block (columns A:H) by one more
# week column (I), mirroring the existing H column's values and formatting,
# to add GPA-iteration sprint data.
$src = $ws.Range("H27:H50")
$dst = $ws.Range("I27:I50")
$src.Copy($dst)

# The new week's "current" progress counters reset to 0 for the stories
# that were still in-flight (H had partial completion counts here; the
# newly added week starts the next sprint at zero).
$ws.Cells.Item(39, 9).Value2 = 0
$ws.Cells.Item(40, 9).Value2 = 0
$ws.Cells.Item(49, 9).Value2 = 0
$ws.Cells.Item(50, 9).Value2 = 0

# Move the active selection to reflect the new bottom-of-sheet working cell.
$ws.Range("J52").Select()
